$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.343.52"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "3.684.23"
$ws.Range("E3").Value = "  -3.35%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "682.69"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.40"
$ws.Range("E6").Value = "  -5.31%  "
$ws.Range("D7").Value = "3.683.13"
$ws.Range("E7").Value = "  -3.38%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -4.39%  "
$ws.Range("E10").Value = "  -8.07%  "
$ws.Range("E11").Value = "  -5.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -5.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.67"
$ws.Range("E14").Value = "  -6.19%  "
$ws.Range("D15").Value = "4.303.87"
$ws.Range("D16").Value = "3.680.75"
$ws.Range("E16").Value = "  -4.08%  "
$ws.Range("D17").Value = "69.412.95"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.35"
$ws.Range("E19").Value = "  -6.30%  "
$ws.Range("E20").Value = "  -7.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.09"
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.89"
$ws.Range("E22").Value = "  -7.54%  "
$ws.Range("E23").Value = "  -7.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.23"
$ws.Range("E24").Value = "  -4.75%  "
$ws.Range("D25").Value = "3.828.62"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("E26").Value = "  -9.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.52"
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.51"
$ws.Range("E29").Value = "  -8.33%  "
$ws.Range("E30").Value = "  -9.80%  "
$ws.Range("E31").Value = "  -11.05%  "
$ws.Range("E32").Value = "  -7.81%  "
$ws.Range("E33").Value = "  -7.02%  "
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.10"
$ws.Range("E35").Value = "  -6.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "3.652.32"
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.49"
$ws.Range("E38").Value = "  -7.28%  "
$ws.Range("E39").Value = "  +6.30%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.33"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0936"
$ws.Range("E41").Value = "  -7.48%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("E44").Value = "  -6.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "162.55"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.31"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "30.11"
$ws.Range("E47").Value = "  +5.53%  "
$ws.Range("E48").Value = "  -13.45%  "
$ws.Range("E49").Value = "  -8.64%  "
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.12"
$ws.Range("E51").Value = "  -3.31%  "
